$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The previous run only processed the first product (row 2): set its
# image-file path (now that downloads record the local path instead of
# leaving the cell blank).
$ws.Range("C2").Value = "hinh_anh_san_pham\SRM_Neutrogena_Ngua_Mun_100ml.jpg"

# Remaining rows haven't been processed yet this run - clear the leftover
# barcode/name values (and any stray blank image cell) so they can be
# re-populated, while keeping their row formatting intact.
$ws.Range("A3:C6").ClearContents()
$ws.Range("A7:B10").ClearContents()

# Rows 7-10 are now reset to the standard (thick-bottom-border) row
# height of 18.6pt, marked as an explicit/custom height.
$ws.Rows.Item(7).RowHeight = 18.6
$ws.Rows.Item(8).RowHeight = 18.6
$ws.Rows.Item(9).RowHeight = 18.6
$ws.Rows.Item(10).RowHeight = 18.6

# Column C needs to be wider to fit the longer image-path text.
$ws.Columns.Item(3).ColumnWidth = 64.21875

# A few extra scratch columns picked up an explicit width while the
# script was working.
$ws.Columns.Item(4).ColumnWidth = 8.88671875
$ws.Columns.Item(5).ColumnWidth = 8.88671875
$ws.Columns.Item(6).ColumnWidth = 8.88671875

# Leave the selection on the next cell to fill in.
$ws.Range("C3").Select()
